$wb = $excel.ActiveWorkbook

# Column F ("想去人数" / "want-to-go count") updates per sheet, generated from the target diff.
$updates = @{
    "展览" = @{ "2"=675; "7"=599; "9"=342; "10"=37; "11"=6065; "12"=656; "14"=22; "15"=255; "16"=178; "17"=90; "18"=591; "19"=20; "21"=36; "22"=178; "23"=1341; "26"=65; "27"=2087; "28"=197; "30"=364; "31"=50; "32"=3360 }
    "演出" = @{ "6"=100; "10"=655; "20"=360; "22"=4061; "26"=160 }
    "本地生活" = @{ "5"=2607; "8"=1510; "12"=676 }
    "全部类型" = @{ "4"=2607; "7"=1510; "10"=675; "12"=676; "13"=599; "14"=100; "16"=342; "17"=37; "18"=6065; "20"=656; "22"=655; "23"=255; "24"=178; "25"=90; "26"=591; "31"=20; "32"=36; "34"=360; "38"=160; "40"=65; "43"=2087; "46"=197; "48"=364; "49"=3360 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $updates[$sheetName]
    foreach ($rowKey in $rowMap.Keys) {
        $row = [int]$rowKey
        $newValue = $rowMap[$rowKey]
        $ws.Cells.Item($row, 6).Value = $newValue
    }
}

Write-Output "Updated F-column want-to-go counts across sheets."
